$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.447.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.753.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.39%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.752.22"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.46%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.539"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.167"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.38"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.462"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000250"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.389.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.772.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.447.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.65%  "
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +16.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "495.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.729"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000153"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +10.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.74%  "
$ws.Range("E30").Value = "  +2.60%  "
$ws.Range("E31").Value = "  +6.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.908.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.72%  "
$ws.Range("E35").Value = "  +1.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.694.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.31%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.79%  "
$ws.Range("E40").Value = "  +2.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.325"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "438.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.36%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.13%  "
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.807.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "140.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0355"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.03%  "
